$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 607.6531
$ws.Range("I15").Value = 607.6531
$ws.Range("K15").Value = 1822.9593
$ws.Range("M15").Value = -1653.9593
$ws.Range("H17").Value = 3711.4375
$ws.Range("J17").Value = 3018.3572
$ws.Range("L17").Value = 9055.071599999999
$ws.Range("N17").Value = -9391.071599999999
$ws.Range("H32").Value = 1802.6666
$ws.Range("I32").Value = 2499.5
$ws.Range("J32").Value = 1603.5714
$ws.Range("K32").Value = 2499.5
$ws.Range("L32").Value = 1603.5714
$ws.Range("M32").Value = -2173.5
$ws.Range("N32").Value = -2255.5714
$ws.Range("H33").Value = 136.4
$ws.Range("I33").Value = 136.4
$ws.Range("K33").Value = 136.4
$ws.Range("M33").Value = 92.59999999999999
$ws.Range("H64").Value = 3112
$ws.Range("I64").Value = 2853.3333
$ws.Range("K64").Value = 2853.3333
$ws.Range("M64").Value = -2605.3333
$ws.Range("H67").Value = 3112
$ws.Range("I67").Value = 2853.3333
$ws.Range("K67").Value = 2853.3333
$ws.Range("M67").Value = -1995.3333
$ws.Range("H98").Value = 1089.9395
$ws.Range("I98").Value = 902.6667
$ws.Range("J98").Value = 2962.6667
$ws.Range("K98").Value = 902.6667
$ws.Range("L98").Value = 2962.6667
$ws.Range("M98").Value = 595.3333
$ws.Range("N98").Value = -5958.6667
$ws.Range("H107").Value = 1583.909
$ws.Range("I107").Value = 804.2857
$ws.Range("K107").Value = 804.2857
$ws.Range("M107").Value = 1115.7143
$ws.Range("H116").Value = 16754
$ws.Range("J116").Value = 2879.6667
$ws.Range("L116").Value = 2879.6667
$ws.Range("N116").Value = -9763.6667
$ws.Range("H122").Value = 1089.9395
$ws.Range("I122").Value = 902.6667
$ws.Range("J122").Value = 2962.6667
$ws.Range("K122").Value = 2708.0001
$ws.Range("L122").Value = 8888.000100000001
$ws.Range("M122").Value = -258.0001000000002
$ws.Range("N122").Value = -13788.0001
$ws.Range("H132").Value = 1198.9584
$ws.Range("I132").Value = 1121.3529
$ws.Range("J132").Value = 1387.4286
$ws.Range("K132").Value = 3364.0587
$ws.Range("L132").Value = 4162.2858
$ws.Range("M132").Value = -834.0587000000005
$ws.Range("N132").Value = -9222.2858
$ws.Range("H141").Value = 3116567.8
$ws.Range("I141").Value = 4002698.2
$ws.Range("K141").Value = 12008094.6
$ws.Range("M141").Value = -12002914.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4856.25
$ws.Range("I32").Value = 3682.4666
$ws.Range("K32").Value = 3682.4666
$ws.Range("M32").Value = -3395.4666
$ws.Range("H61").Value = 10513.583
$ws.Range("I61").Value = 17473.4
$ws.Range("K61").Value = 17473.4
$ws.Range("M61").Value = -17261.4
$ws.Range("H74").Value = 1291.174
$ws.Range("I74").Value = 500
$ws.Range("K74").Value = 500
$ws.Range("M74").Value = 374
$ws.Range("H77").Value = 1291.174
$ws.Range("I77").Value = 500
$ws.Range("K77").Value = 2500
$ws.Range("M77").Value = 1868
$ws.Range("H97").Value = 1119.3
$ws.Range("I97").Value = 1065.9445
$ws.Range("K97").Value = 1065.9445
$ws.Range("M97").Value = -569.9445000000001
$ws.Range("H102").Value = 1451.4445
$ws.Range("I102").Value = 1328
$ws.Range("J102").Value = 1513.1666
$ws.Range("K102").Value = 1328
$ws.Range("L102").Value = 1513.1666
$ws.Range("M102").Value = 294
$ws.Range("N102").Value = -4757.1666
$ws.Range("H110").Value = 1393.4445
$ws.Range("I110").Value = 284.42856
$ws.Range("K110").Value = 284.42856
$ws.Range("M110").Value = 1760.57144
$ws.Range("H132").Value = 2733.8462
$ws.Range("I132").Value = 2218.4285
$ws.Range("J132").Value = 3335.1667
$ws.Range("K132").Value = 6655.2855
$ws.Range("L132").Value = 10005.5001
$ws.Range("M132").Value = -4125.2855
$ws.Range("N132").Value = -15065.5001
$ws.Range("H136").Value = 10513.583
$ws.Range("I136").Value = 17473.4
$ws.Range("K136").Value = 52420.2
$ws.Range("M136").Value = -49870.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7825.4443
$ws.Range("I134").Value = 9052.5
$ws.Range("K134").Value = 27157.5
$ws.Range("M134").Value = -24622.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6214110
$ws.Range("J58").Value = 4500
$ws.Range("L58").Value = 4500
$ws.Range("N58").Value = -4906
$ws.Range("H99").Value = 1251808.6
$ws.Range("I99").Value = 2501742.2
$ws.Range("J99").Value = 1875
$ws.Range("K99").Value = 2501742.2
$ws.Range("L99").Value = 1875
$ws.Range("M99").Value = -2500244.2
$ws.Range("N99").Value = -4871
$ws.Range("H126").Value = 1251808.6
$ws.Range("I126").Value = 2501742.2
$ws.Range("J126").Value = 1875
$ws.Range("K126").Value = 7505226.600000001
$ws.Range("L126").Value = 5625
$ws.Range("M126").Value = -7502756.600000001
$ws.Range("N126").Value = -10565
$ws.Range("H132").Value = 2423.8262
$ws.Range("I132").Value = 1572.9412
$ws.Range("K132").Value = 4718.8236
$ws.Range("M132").Value = -2188.8236
$ws.Range("H136").Value = 6214110
$ws.Range("J136").Value = 4500
$ws.Range("L136").Value = 13500
$ws.Range("N136").Value = -18600
$ws.Range("H141").Value = 66599.8
$ws.Range("J141").Value = 65249.75
$ws.Range("L141").Value = 65249.75
$ws.Range("N141").Value = -75609.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 938.5
$ws.Range("J80").Value = 751.3333
$ws.Range("L80").Value = 2253.9999
$ws.Range("N80").Value = -4125.9999
$ws.Range("H83").Value = 938.5
$ws.Range("J83").Value = 751.3333
$ws.Range("L83").Value = 6761.9997
$ws.Range("N83").Value = -16121.9997

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2268.818
$ws.Range("I102").Value = 1823.625
$ws.Range("K102").Value = 1823.625
$ws.Range("M102").Value = -201.625

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10398.889
$ws.Range("I40").Value = 9770.429
$ws.Range("J40").Value = 12598.5
$ws.Range("K40").Value = 9770.429
$ws.Range("L40").Value = 12598.5
$ws.Range("M40").Value = -9634.429
$ws.Range("N40").Value = -12870.5
$ws.Range("H46").Value = 1963.25
$ws.Range("J46").Value = 3155.8
$ws.Range("L46").Value = 3155.8
$ws.Range("N46").Value = -3531.8
$ws.Range("H55").Value = 582.5263
$ws.Range("J55").Value = 562.1818
$ws.Range("L55").Value = 562.1818
$ws.Range("N55").Value = -908.1818
$ws.Range("H61").Value = 2438.1904
$ws.Range("I61").Value = 2430.923
$ws.Range("J61").Value = 2450
$ws.Range("K61").Value = 2430.923
$ws.Range("L61").Value = 2450
$ws.Range("M61").Value = -2228.923
$ws.Range("N61").Value = -2854
$ws.Range("I93").Value = 869
$ws.Range("K93").Value = 869
$ws.Range("M93").Value = 379
$ws.Range("H113").Value = 2438.1904
$ws.Range("I113").Value = 2430.923
$ws.Range("J113").Value = 2450
$ws.Range("K113").Value = 2430.923
$ws.Range("L113").Value = 2450
$ws.Range("M113").Value = -260.9229999999998
$ws.Range("N113").Value = -6790
$ws.Range("H132").Value = 1968.8334
$ws.Range("I132").Value = 1419.35
$ws.Range("J132").Value = 2655.6875
$ws.Range("K132").Value = 4258.049999999999
$ws.Range("L132").Value = 7967.0625
$ws.Range("M132").Value = -1728.049999999999
$ws.Range("N132").Value = -13027.0625
$ws.Range("H136").Value = 3300.7144
$ws.Range("I136").Value = 3767
$ws.Range("K136").Value = 11301
$ws.Range("M136").Value = -8751

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1112.25
$ws.Range("I81").Value = 1333
$ws.Range("K81").Value = 2666
$ws.Range("M81").Value = -1605
$ws.Range("H84").Value = 1112.25
$ws.Range("I84").Value = 1333
$ws.Range("K84").Value = 13330
$ws.Range("M84").Value = -8026
$ws.Range("H132").Value = 1364.6207
$ws.Range("I132").Value = 968.6585
$ws.Range("J132").Value = 2319.5881
$ws.Range("K132").Value = 2905.9755
$ws.Range("L132").Value = 6958.7643
$ws.Range("M132").Value = -375.9755
$ws.Range("N132").Value = -12018.7643
